$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update revised values in rows 174-176 (monthly series were restated) ---
$ws.Range("N174").Value = 15220
$ws.Range("O174").Value = 17228
$ws.Range("P174").Value = 7333
$ws.Range("Q174").Value = 8940

$ws.Range("M175").Value = 7393
$ws.Range("N175").Value = 14654
$ws.Range("O175").Value = 17854
$ws.Range("P175").Value = 7601
$ws.Range("Q175").Value = 8509

$ws.Range("L176").Value = 838
$ws.Range("M176").Value = 5694
$ws.Range("N176").Value = 13281
$ws.Range("O176").Value = 17297
$ws.Range("P176").Value = 8738
$ws.Range("Q176").Value = 8331

# --- Append new row 177 for period 01-08-2021 ---
# Column A holds a date-looking label that must stay a plain text string
# (matching the rest of column A), so we stage it as a formula result
# (guaranteed text) in a scratch cell, copy/paste the value across, and
# then wipe the scratch cell completely -- this avoids Excel's automatic
# date recognition (which would otherwise turn the text into a serial
# date number) while also avoiding introducing any new cell style.
$staging = $ws.Range("ZZ1")
$staging.Formula = '="01-08-2021"'
$staging.Copy()
$ws.Range("A177").PasteSpecial(-4163)
$staging.Clear()

$ws.Range("B177").Value = -21728
$ws.Range("C177").Value = -14674
$ws.Range("D177").Value = 18271
$ws.Range("E177").Value = 32945
$ws.Range("F177").Value = -7054
$ws.Range("G177").Value = 4323
$ws.Range("H177").Value = 11377
$ws.Range("I177").Value = 66916
$ws.Range("J177").Value = 51216
$ws.Range("K177").Value = 371
$ws.Range("L177").Value = 2360
$ws.Range("M177").Value = 5337
$ws.Range("N177").Value = 10517
$ws.Range("O177").Value = 16272
$ws.Range("P177").Value = 8136
$ws.Range("Q177").Value = 8223
$ws.Range("R177").Value = 15700
